$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "D" (Price) column cells hold plain-looking numeric text (e.g. "588.17",
# "1.00", "70.467.43") that must stay literal text, matching the workbook's
# original inlineStr cell type. Force a Text number format before assigning so
# Excel does not silently coerce the string into a float and mangle it
# (dropping trailing zeros, losing thousand-dot grouping, introducing FP noise).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.467.43"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.559.97"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.17"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.45"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.552.94"
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.217"
$ws.Range("E10").Value = "  +11.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.15"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.42"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.124.33"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.430.75"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.573.96"
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.77"
$ws.Range("E18").Value = "  +5.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.99"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "581.67"
$ws.Range("E20").Value = "  +9.38%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.994"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.91"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.63"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.86"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.76"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  +4.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.36"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.20"
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.54"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.76"
$ws.Range("E35").Value = "  +28.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.23"
$ws.Range("E36").Value = "  +6.73%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.713.44"
$ws.Range("E37").Value = "  +11.71%  "
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "526.12"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.88"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("E43").Value = "  +6.34%  "
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0455"
$ws.Range("E45").Value = "  +5.82%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.22"
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.42"
$ws.Range("E51").Value = "  +8.90%  "
